$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 ---------------------------------------------------------------
# Copy formatting from the row above (row 6) / header rows so the new cells
# reuse the existing cell styles (date, time, etc.) instead of creating new
# duplicate style entries.
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = 40220

$ws.Range("B6").Copy($ws.Range("B7"))
$ws.Range("B7").Value = 0.4375

$ws.Range("C6").Copy($ws.Range("C7"))
$ws.Range("C7").Value = 0.60416666666666663

$ws.Range("D2").Copy($ws.Range("D7"))
$ws.Range("D7").Value = 4

$ws.Range("E4").Copy($ws.Range("E7"))
$ws.Range("E7").Value = "Design"

$ws.Range("F2").Copy($ws.Range("F7"))
$ws.Range("F7").Value = "MessageQueue and rendering design, Cutting out tiles"

# --- Row 8 ---------------------------------------------------------------
$ws.Range("A6").Copy($ws.Range("A8"))
$ws.Range("A8").Value = 40225

$ws.Range("B6").Copy($ws.Range("B8"))
$ws.Range("B8").Value = 0.45833333333333331

$ws.Range("C6").Copy($ws.Range("C8"))
$ws.Range("C8").Value = 0.60416666666666663

# D8 holds the text "3.5" (not the number 3.5). Force text entry with a
# leading apostrophe, then reset the cell's style to the plain default style
# (borrowed from a cell that already carries that style) so no extra
# "quote prefix" style gets created.
$ws.Range("D8").Value = "'3.5"
$ws.Range("D8").Style = $ws.Range("D2").Style

$ws.Range("E6").Copy($ws.Range("E8"))
$ws.Range("E8").Value = "Analysis"

$ws.Range("F6").Copy($ws.Range("F8"))
$ws.Range("F8").Value = "Plan of Attack"

# --- Selection -------------------------------------------------------------
$ws.Range("F17").Select()
